$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.183.02"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "1.854.09"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.51"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6859"
$ws.Range("E6").Value = "  -4.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07675"
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3047"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.17"
$ws.Range("E10").Value = "  -5.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08109"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.881.90"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7240"
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.191"
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.53"
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("D16").Value = "29.188.17"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007811"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.726"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.23"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "234.12"
$ws.Range("E20").Value = "  -4.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "2.102.08"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.459"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.70"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.941"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1429"
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.959"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.505"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.005"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05161"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.184"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7059"
$ws.Range("E36").Value = "  -3.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.028"
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01849"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.679"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9153"
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("D42").Value = "1.103.61"
$ws.Range("E42").Value = "  +6.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.958"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4283"
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.14"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.75"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.781"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").Value = "2.000.74"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.149"
$ws.Range("E50").Value = "  -5.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.988"
$ws.Range("E51").Value = "  -6.22%  "
